# Jogos da Semana FlashScore 2024-10-31 update
# - Remove the old row 5 (Al Shabab vs Al Wehda, 11:50) -- all rows below shift up one
# - The row that shifts into row 6 (formerly row 7: Al Ittihad vs Al Ahli SC) gets a couple
#   of cell corrections (Odd_Over05_FT/Odd_Under05_FT order fix, Odd_CS_3-3_HT value)
# - Two brand-new fixtures are appended as rows 7 and 8 (Switzerland Super League)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete first fixture row; everything below shifts up automatically.
$ws.Rows(5).Delete()

$row6Data = @(
    "G2TuVbho",
    "31/10/2024",
    "15:00",
    "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE",
    "Al Ittihad",
    "Al Ahli SC",
    2.2,
    3.7,
    2.8,
    2.63,
    2.5,
    3,
    1.02,
    11,
    1.13,
    5.5,
    1.44,
    2.63,
    1.22,
    4,
    1.4,
    2.75,
    15,
    15,
    10,
    23,
    15,
    19,
    21,
    8,
    11,
    29,
    81,
    17,
    19,
    12,
    29,
    19,
    21,
    5,
    11,
    17,
    34,
    41,
    81,
    4,
    7,
    41,
    5.5,
    13,
    17,
    41,
    51,
    81,
    350,
    81
)

$row7Data = @(
    "Q5UoY3ef",
    "31/10/2024",
    "16:30",
    "SWITZERLAND - SUPER LEAGUE",
    "Grasshoppers",
    "Lugano",
    3.4,
    3.4,
    2.05,
    3.75,
    2.25,
    2.75,
    1.04,
    12,
    1.25,
    3.75,
    1.83,
    2.03,
    1.36,
    3,
    1.67,
    2.1,
    12,
    17,
    12,
    34,
    26,
    29,
    12,
    6.5,
    13,
    41,
    151,
    8.5,
    11,
    9,
    19,
    17,
    23,
    5.5,
    17,
    23,
    51,
    67,
    151,
    3,
    7.5,
    51,
    4.33,
    11,
    21,
    41,
    51,
    126,
    81,
    81
)

$row8Data = @(
    "YVXgWsQ6",
    "31/10/2024",
    "16:30",
    "SWITZERLAND - SUPER LEAGUE",
    "Servette",
    "Luzern",
    1.81,
    3.7,
    3.9,
    2.4,
    2.4,
    4,
    1.03,
    17,
    1.17,
    5,
    1.57,
    2.35,
    1.29,
    3.5,
    1.53,
    2.38,
    10,
    11,
    8.5,
    17,
    13,
    21,
    15,
    7.5,
    12,
    34,
    126,
    15,
    23,
    13,
    41,
    26,
    29,
    4.33,
    9.5,
    17,
    29,
    41,
    101,
    3.5,
    7,
    41,
    6,
    19,
    23,
    51,
    67,
    126,
    351,
    81
)

for ($i = 0; $i -lt $row6Data.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6Data[$i]
}

for ($i = 0; $i -lt $row7Data.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $row7Data[$i]
}

for ($i = 0; $i -lt $row8Data.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $row8Data[$i]
}

Write-Output "Update complete"

